# corrected data source reset at startTest.
# Replace the old sample (dates/numbers) with a fresh 4-column x 3-row
# block of plain text values coming from the (reset) data source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe whatever was there before (values + number formats on A1:C3).
$ws.Cells.Clear()

# New data: columns A-D, rows 1-3 -> "<Col> 1<Row>"
$data = @(
    @("A 11", "B 11", "C 11", "D 11"),
    @("A 12", "B 12", "C 12", "D 12"),
    @("A 13", "B 13", "C 13", "D 13")
)

for ($r = 1; $r -le 3; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}

# Column width adjustments (column C slightly narrower, new column D sized
# to fit its content).
$ws.Columns.Item(3).ColumnWidth = 13.3
$ws.Columns.Item(4).ColumnWidth = 10.084

# Move/restore the active selection like the author's session.
$ws.Range("C9").Select()
